# SLSkeywordsremoveandadd/url.xlsx
#
# The "LSS NT" keyword column (E) is removed entirely, and four of the
# remaining keyword columns get their text replaced with shorter/renamed
# keywords:
#   G: "NT SCIENCE"             -> "n(t)"
#   H: "NT Science"              -> "NT sci"
#   J: "Lower Secondary NT"      -> "nts"
#   K: "Normal Technical Science"-> "n(t)"
# Columns D, F, I, L, M keep their existing text (lss nt / nt science /
# G1 science / G1 LSS / G1 Science).
#
# The data rows run from row 2 to row 50.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the keyword columns in the order the author appears to have typed
# them (this governs the order new strings are appended to the shared
# string table): J, then G, then H, then K.
$ws.Range("J2:J50").Value = "nts"
$ws.Range("G2:G50").Value = "n(t)"
$ws.Range("H2:H50").Value = "NT sci"
$ws.Range("K2:K50").Value = "n(t)"

# Remove the "LSS NT" keyword column entirely (E2:E50).
$ws.Range("E2:E50").ClearContents()

# Match the author's final selection/scroll state recorded in the sheet
# view (column E was being worked on before being cleared).
$ws.Range("E2:E50").Select()
